# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 220
    $ws.Range("F3").Value = 257
    $ws.Range("F7").Value = 6265
    $ws.Range("F15").Value = 463
}
